# Fruta / hortaliza, semanal
#
# The weekly refresh prepends a new date block (date serial 44641) of three
# rows ("Pintón" / "Primera Maduro" / "Primera Pintón") right before the
# existing row 494, pushing all the following data rows down by three rows
# (494-583 -> 497-586). The sheet therefore grows from A1:T583 to A1:T586.
#
# All of columns A, B, C, E, F, G, H, I, J, K, Q, R, T are constant for this
# whole "Terminal La Palmera de La Serena / Plátano" block, and column L
# cycles with period 3 ("Pintón", "Primera Maduro", "Primera Pintón"), so the
# freshly inserted rows simply reuse that same template; only D (fecha), M
# (volumen, kept the same as the pattern), N/O/P (precios) and S (precio
# $/kg) carry new values for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 494, shifting the existing 494:583 block down to
# 497:586 (this also grows dimension/UsedRange to A1:T586 automatically).
$ws.Rows("494:496").Insert()

# --- Row 494: Pintón, volumen 80 ---
$ws.Cells.Item(494, 1).Value = 8
$ws.Cells.Item(494, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(494, 3).Value = "Coquimbo"
$ws.Cells.Item(494, 4).Value = 44641
$ws.Cells.Item(494, 5).Value = 4
$ws.Cells.Item(494, 6).Value = "Fruta"
$ws.Cells.Item(494, 7).Value = 100108
$ws.Cells.Item(494, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(494, 9).Value = 100108006
$ws.Cells.Item(494, 10).Value = "Plátano"
$ws.Cells.Item(494, 11).Value = "Sin especificar"
$ws.Cells.Item(494, 12).Value = "Pintón"
$ws.Cells.Item(494, 13).Value = 80
$ws.Cells.Item(494, 14).Value = 18000
$ws.Cells.Item(494, 15).Value = 18000
$ws.Cells.Item(494, 16).Value = 18000
$ws.Cells.Item(494, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(494, 18).Value = "Ecuador"
$ws.Cells.Item(494, 19).Value = 900
$ws.Cells.Item(494, 20).Value = 20

# --- Row 495: Primera Maduro, volumen 120 ---
$ws.Cells.Item(495, 1).Value = 8
$ws.Cells.Item(495, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(495, 3).Value = "Coquimbo"
$ws.Cells.Item(495, 4).Value = 44641
$ws.Cells.Item(495, 5).Value = 4
$ws.Cells.Item(495, 6).Value = "Fruta"
$ws.Cells.Item(495, 7).Value = 100108
$ws.Cells.Item(495, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(495, 9).Value = 100108006
$ws.Cells.Item(495, 10).Value = "Plátano"
$ws.Cells.Item(495, 11).Value = "Sin especificar"
$ws.Cells.Item(495, 12).Value = "Primera Maduro"
$ws.Cells.Item(495, 13).Value = 120
$ws.Cells.Item(495, 14).Value = 20000
$ws.Cells.Item(495, 15).Value = 20000
$ws.Cells.Item(495, 16).Value = 20000
$ws.Cells.Item(495, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(495, 18).Value = "Ecuador"
$ws.Cells.Item(495, 19).Value = 1000
$ws.Cells.Item(495, 20).Value = 20

# --- Row 496: Primera Pintón, volumen 120 ---
$ws.Cells.Item(496, 1).Value = 8
$ws.Cells.Item(496, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(496, 3).Value = "Coquimbo"
$ws.Cells.Item(496, 4).Value = 44641
$ws.Cells.Item(496, 5).Value = 4
$ws.Cells.Item(496, 6).Value = "Fruta"
$ws.Cells.Item(496, 7).Value = 100108
$ws.Cells.Item(496, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(496, 9).Value = 100108006
$ws.Cells.Item(496, 10).Value = "Plátano"
$ws.Cells.Item(496, 11).Value = "Sin especificar"
$ws.Cells.Item(496, 12).Value = "Primera Pintón"
$ws.Cells.Item(496, 13).Value = 120
$ws.Cells.Item(496, 14).Value = 21000
$ws.Cells.Item(496, 15).Value = 21000
$ws.Cells.Item(496, 16).Value = 21000
$ws.Cells.Item(496, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(496, 18).Value = "Ecuador"
$ws.Cells.Item(496, 19).Value = 1050
$ws.Cells.Item(496, 20).Value = 20
